$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit re-shuffles the per-observation data across rows 4-9 (a single
# species-records columns get redistributed among the existing rows), while
# row-invariant / site-level columns are left untouched.
#
# Mapping: new content of row R = old content of row Src
#   4 <- 8
#   5 <- 4
#   6 <- 9
#   7 <- 5
#   8 <- 6
#   9 <- 7
$mapping = @{
    4 = 8
    5 = 4
    6 = 9
    7 = 5
    8 = 6
    9 = 7
}

# Columns that move together with each observation's row.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "W", "Z", "AB")

# Snapshot current (pre-edit) values for every affected cell before writing
# anything, since the remap is a single 6-cycle and would otherwise clobber
# source data mid-way.
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    foreach ($col in $cols) {
        $key = "$col$row"
        $snapshot[$key] = $ws.Range($key).Value2
    }
}

# Write the remapped values.
foreach ($row in $mapping.Keys) {
    $src = $mapping[$row]
    foreach ($col in $cols) {
        $srcKey = "$col$src"
        $dstKey = "$col$row"
        $ws.Range($dstKey).Value2 = $snapshot[$srcKey]
    }
}
